$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = '45765771'
$ws.Range("D16").Value = 'ELIFELET ROCHA GARCIA'
$ws.Range("E16").Value = '2304'
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1300000

$ws.Range("C17").Value = '1047375750'
$ws.Range("D17").Value = 'JOHANIS IBAÑEZ'
$ws.Range("E17").Value = '2304'
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = '1143340401'
$ws.Range("D18").Value = 'LEIDIS MARCELA PEÑA ROBLES'
$ws.Range("E18").Value = '2304'
$ws.Range("F18").Value = 46400
$ws.Range("G18").Value = 1300000

$ws.Range("C19").Value = '49759568'
$ws.Range("D19").Value = 'LILIANA ACUÑA BRAVO'
$ws.Range("E19").Value = '2305'
$ws.Range("F19").Value = 46400
$ws.Range("G19").Value = 1300000

$ws.Range("C20").Value = '49759568'
$ws.Range("D20").Value = 'LILIANA ACUÑA BRAVO'
$ws.Range("E20").Value = '2304'
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1300000

$ws.Range("C21").Value = '73571102'
$ws.Range("D21").Value = 'RUBEN DARIO VILLAMIL AGUILAR'
$ws.Range("E21").Value = '2305'
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1300000

$ws.Range("C22").Value = '73571102'
$ws.Range("D22").Value = 'RUBEN DARIO VILLAMIL AGUILAR'
$ws.Range("E22").Value = '2304'
$ws.Range("F22").Value = 41760
$ws.Range("G22").Value = 1300000

$ws.Range("C23").Value = '45523428'
$ws.Range("D23").Value = 'KATHERINE GISEL BRITO SERRANO'
$ws.Range("E23").Value = '2305'
$ws.Range("F23").Value = 46400
$ws.Range("G23").Value = 1300000

$ws.Range("C24").Value = '45523428'
$ws.Range("D24").Value = 'KATHERINE GISEL BRITO SERRANO'
$ws.Range("E24").Value = '2304'
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1300000

$ws.Range("C25").Value = '45477182'
$ws.Range("D25").Value = 'YANERYS CASTILLO ORTEGA'
$ws.Range("E25").Value = '2305'
$ws.Range("F25").Value = 46400
$ws.Range("G25").Value = 1160000

$ws.Range("C26").Value = '45477182'
$ws.Range("D26").Value = 'YANERYS CASTILLO ORTEGA'
$ws.Range("E26").Value = '2304'
$ws.Range("F26").Value = 46400
$ws.Range("G26").Value = 1160000

$ws.Range("C27").Value = '45525282'
$ws.Range("D27").Value = 'YESMITH DEL CARMEN MARRUGO ORTIZ'
$ws.Range("E27").Value = '2304'
$ws.Range("F27").Value = 46400
$ws.Range("G27").Value = 1160000

$ws.Range("C28").Value = '45491200'
$ws.Range("D28").Value = 'NORIS DE JESUS TORRES QUINTANA'
$ws.Range("E28").Value = '2305'
$ws.Range("F28").Value = 53568
$ws.Range("G28").Value = 1473120

$ws.Range("C29").Value = '45491200'
$ws.Range("D29").Value = 'NORIS DE JESUS TORRES QUINTANA'
$ws.Range("E29").Value = '2304'
$ws.Range("F29").Value = 53568
$ws.Range("G29").Value = 1473120

$ws.Range("C30").Value = '1143326275'
$ws.Range("D30").Value = 'MARIA EUGENIA DIAZ RUIZ'
$ws.Range("E30").Value = '2304'
$ws.Range("F30").Value = 46400
$ws.Range("G30").Value = 1160000

$ws.Range("C31").Value = '45554005'
$ws.Range("D31").Value = 'MARTHA FAJARDO CHICO'
$ws.Range("E31").Value = '2305'
$ws.Range("F31").Value = 49600
$ws.Range("G31").Value = 1390000

$ws.Range("C32").Value = '45554005'
$ws.Range("D32").Value = 'MARTHA FAJARDO CHICO'
$ws.Range("E32").Value = '2304'
$ws.Range("F32").Value = 49600
$ws.Range("G32").Value = 1390000

$ws.Range("C33").Value = '15025168'
$ws.Range("D33").Value = 'BLAS JOSE HERNANDEZ GENES'
$ws.Range("E33").Value = '2309'
$ws.Range("F33").Value = 46400
$ws.Range("G33").Value = 1160000

$ws.Range("C34").Value = '1050949359'
$ws.Range("D34").Value = 'INGRID CANTILLO ROCHA'
$ws.Range("E34").Value = '2304'
$ws.Range("F34").Value = 46400
$ws.Range("G34").Value = 1300000

$ws.Range("C35").Value = '1070814213'
$ws.Range("D35").Value = 'ERICA ROCIO GONZALEZ VARGAS'
$ws.Range("E35").Value = '2304'
$ws.Range("F35").Value = 46400
$ws.Range("G35").Value = 1300000

$ws.Range("C36").Value = '26795657'
$ws.Range("D36").Value = 'BEATRIZ ELENA NOVOA TAPIA'
$ws.Range("E36").Value = '2305'
$ws.Range("F36").Value = 46400
$ws.Range("G36").Value = 1300000

$ws.Range("C37").Value = '26795657'
$ws.Range("D37").Value = 'BEATRIZ ELENA NOVOA TAPIA'
$ws.Range("E37").Value = '2304'
$ws.Range("F37").Value = 46400
$ws.Range("G37").Value = 1300000

$ws.Range("C38").Value = '22464424'
$ws.Range("D38").Value = 'IRINA CHARRIS HOYOS'
$ws.Range("E38").Value = '2305'
$ws.Range("F38").Value = 46400
$ws.Range("G38").Value = 1160000

$ws.Range("C39").Value = '22464424'
$ws.Range("D39").Value = 'IRINA CHARRIS HOYOS'
$ws.Range("E39").Value = '2304'
$ws.Range("F39").Value = 46400
$ws.Range("G39").Value = 1160000

$ws.Range("C40").Value = '1047409863'
$ws.Range("D40").Value = 'ERIKA ROSA BERRIO BARRIOS'
$ws.Range("E40").Value = '2304'
$ws.Range("F40").Value = 46400
$ws.Range("G40").Value = 1160000

$ws.Range("C41").Value = '73154999'
$ws.Range("D41").Value = 'FELIX ANTONIO BELLIDO MONTERO'
$ws.Range("E41").Value = '2305'
$ws.Range("F41").Value = 46400
$ws.Range("G41").Value = 1300000

$ws.Range("C42").Value = '73154999'
$ws.Range("D42").Value = 'FELIX ANTONIO BELLIDO MONTERO'
$ws.Range("E42").Value = '2304'
$ws.Range("F42").Value = 46400
$ws.Range("G42").Value = 1300000

$ws.Range("C43").Value = '33102950'
$ws.Range("D43").Value = 'LUZDARY HUERTAS ARANZA'
$ws.Range("E43").Value = '2305'
$ws.Range("F43").Value = 46400
$ws.Range("G43").Value = 1300000

$ws.Range("C44").Value = '33102950'
$ws.Range("D44").Value = 'LUZDARY HUERTAS ARANZA'
$ws.Range("E44").Value = '2304'
$ws.Range("F44").Value = 46400
$ws.Range("G44").Value = 1300000

$ws.Range("C45").Value = '1047456149'
$ws.Range("D45").Value = 'GUSTAVO ADOLFO JAM FIGUEROA'
$ws.Range("E45").Value = '2304'
$ws.Range("F45").Value = 46400
$ws.Range("G45").Value = 1300000

$ws.Range("C46").Value = '1128058925'
$ws.Range("D46").Value = 'YESENIA PARRA MARIMON'
$ws.Range("E46").Value = '2304'
$ws.Range("F46").Value = 46400
$ws.Range("G46").Value = 1300000

$ws.Range("C47").Value = '1047504667'
$ws.Range("D47").Value = 'WENDY PAOLA GODOY MEDRANO'
$ws.Range("E47").Value = '2306'
$ws.Range("F47").Value = 46400
$ws.Range("G47").Value = 1300000

$ws.Range("C48").Value = '1047504667'
$ws.Range("D48").Value = 'WENDY PAOLA GODOY MEDRANO'
$ws.Range("E48").Value = '2305'
$ws.Range("F48").Value = 46400
$ws.Range("G48").Value = 1300000

$ws.Range("C49").Value = '1047504667'
$ws.Range("D49").Value = 'WENDY PAOLA GODOY MEDRANO'
$ws.Range("E49").Value = '2304'
$ws.Range("F49").Value = 46400
$ws.Range("G49").Value = 1300000

$ws.Range("C50").Value = '1047455773'
$ws.Range("D50").Value = 'EDRIN ALEJANDRO VILLAR VASQUEZ'
$ws.Range("E50").Value = '2304'
$ws.Range("F50").Value = 46400
$ws.Range("G50").Value = 1300000

$ws.Range("C51").Value = '1007314298'
$ws.Range("D51").Value = 'MILEIDIS DEL CARMEN PALACIN VANEGAS'
$ws.Range("E51").Value = '2304'
$ws.Range("F51").Value = 46400
$ws.Range("G51").Value = 1160000
